$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new column F header + values for guest account support
$ws.Range("F1").Value = "is_authenticated"
$ws.Range("F2").Value = 1
$ws.Range("F3").Value = 1
$ws.Range("F4").Value = 1

# Update selection to match target state
$ws.Range("F5").Select()
